# Insert a new weekly price record at row 166 (Arveja Verde, Vega Central
# Mapocho de Santiago), shifting the existing rows 166-175 down to 167-176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 166..175 down one row, creating a blank row 166.
$ws.Rows.Item(166).Insert()

# Populate the new row 166 with the new record.
$ws.Range("A166").Value = 9
$ws.Range("B166").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C166").Value = "Metropolitana"
$ws.Range("D166").Value = 45132
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = 100112022
$ws.Range("G166").Value = "Arveja Verde"
$ws.Range("H166").Value = "Perfection"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 52
$ws.Range("K166").Value = 26000
$ws.Range("L166").Value = 28000
$ws.Range("M166").Value = 27000
$ws.Range("N166").Value = "$/malla 25 kilos"
$ws.Range("O166").Value = "Provincia de Limarí"
$ws.Range("P166").Value = 1080
$ws.Range("Q166").Value = 25
$ws.Range("R166").Value = "Hortaliza"
